# Scheduled-runner update: refresh market-price columns (H-N) on the
# Pandaemonium_Profits leve-crafting sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR)
# with freshly pulled Universalis price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 3: One for the Books / Leather Grimoire
$ws.Range("H3").Value = 50000
$ws.Range("J3").Value = 50000
$ws.Range("L3").Value = 50000
$ws.Range("N3").Value = -50228

# Row 17: One for the Road / Potion
$ws.Range("H17").Value = 929.2593000000001
$ws.Range("J17").Value = 929.2593000000001
$ws.Range("L17").Value = 2787.7779
$ws.Range("N17").Value = -3123.7779

# Row 43: Growing Is Knowing / Growth Formula Gamma
$ws.Range("H43").Value = 1791.8
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 1791.8
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 1791.8
$ws.Range("N43").Value = -1929.8
$ws.Range("M43").Value = ""

# Row 58: A Matter of Vital Importance / Mega-Potion of Vitality
$ws.Range("H58").Value = 112784.89
$ws.Range("J58").Value = 202979.8
$ws.Range("L58").Value = 608939.3999999999
$ws.Range("N58").Value = -609239.3999999999

# Row 70: Consecrating Congregation / Holy Water
$ws.Range("H70").Value = 1561.75
$ws.Range("I70").Value = 1187.3334
$ws.Range("J70").Value = 2043.1428
$ws.Range("K70").Value = 3562.0002
$ws.Range("L70").Value = 6129.428400000001
$ws.Range("M70").Value = -3292.0002
$ws.Range("N70").Value = -6669.428400000001

# Row 73: Curbing the Contagion (L) / Holy Water
$ws.Range("H73").Value = 1561.75
$ws.Range("I73").Value = 1187.3334
$ws.Range("J73").Value = 2043.1428
$ws.Range("K73").Value = 3562.0002
$ws.Range("L73").Value = 6129.428400000001
$ws.Range("M73").Value = -2626.0002
$ws.Range("N73").Value = -8001.428400000001

# Row 74: Adhesive of Antipathy / Wing Glue
$ws.Range("H74").Value = 4063.6365
$ws.Range("I74").Value = 3650
$ws.Range("J74").Value = 5166.6665
$ws.Range("K74").Value = 3650
$ws.Range("L74").Value = 5166.6665
$ws.Range("M74").Value = -2714
$ws.Range("N74").Value = -7038.6665

# Row 77: It's Gonna Grow Back (L) / Wing Glue
$ws.Range("H77").Value = 4063.6365
$ws.Range("I77").Value = 3650
$ws.Range("J77").Value = 5166.6665
$ws.Range("K77").Value = 18250
$ws.Range("L77").Value = 25833.3325
$ws.Range("M77").Value = -13570
$ws.Range("N77").Value = -35193.3325

# Row 87: There Was a Late Fee / Noble Gold
$ws.Range("H87").Value = 39418
$ws.Range("J87").Value = 39418
$ws.Range("L87").Value = 39418
$ws.Range("N87").Value = -41914

# Row 90: A Gate Arcane Is Dragon's Bane (L) / Noble Gold
$ws.Range("H90").Value = 39418
$ws.Range("J90").Value = 39418
$ws.Range("L90").Value = 118254
$ws.Range("N90").Value = -130734

# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 2378.244
$ws.Range("I98").Value = 2378.244
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 2378.244
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -880.2440000000001
$ws.Range("N98").Value = ""

# Row 102: Spell-rebound / Marid Leather Grimoire
$ws.Range("H102").Value = 50000
$ws.Range("J102").Value = 50000
$ws.Range("L102").Value = 50000
$ws.Range("N102").Value = -56490

# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 2378.244
$ws.Range("I122").Value = 2378.244
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7134.732
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4684.732
$ws.Range("N122").Value = ""

# Row 135: For Tired Minds / Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 88238050
$ws.Range("I135").Value = 41669610
$ws.Range("J135").Value = 200002300
$ws.Range("K135").Value = 375026490
$ws.Range("L135").Value = 1800020700
$ws.Range("M135").Value = -375023955
$ws.Range("N135").Value = -1800025770

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 1061101.8
$ws.Range("I138").Value = 1635.9474
$ws.Range("K138").Value = 4907.8422
$ws.Range("M138").Value = 232.1578

# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 3359.923
$ws.Range("I141").Value = 3788.1
$ws.Range("J141").Value = 1932.6666
$ws.Range("K141").Value = 11364.3
$ws.Range("L141").Value = 5797.9998
$ws.Range("M141").Value = -6184.299999999999
$ws.Range("N141").Value = -16157.9998

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 11528.154
$ws.Range("I32").Value = 8495.394
$ws.Range("J32").Value = 28208.334
$ws.Range("K32").Value = 8495.394
$ws.Range("L32").Value = 28208.334
$ws.Range("M32").Value = -8208.394
$ws.Range("N32").Value = -28782.334

# Row 97: Ore for Me / High Steel Ingot
$ws.Range("H97").Value = 2451.5
$ws.Range("I97").Value = 1404.5
$ws.Range("K97").Value = 1404.5
$ws.Range("M97").Value = -908.5

# Row 125: The Incomplete Costume / High Durium Armor of Fending
$ws.Range("H125").Value = 65932.914
$ws.Range("J125").Value = 65932.914
$ws.Range("L125").Value = 65932.914
$ws.Range("N125").Value = -75772.914

$ws = $wb.Worksheets.Item("BSM")
# Row 25: Tools of the Trade / Iron Doming Hammer
$ws.Range("H25").Value = 2701.2856
$ws.Range("I25").Value = 1651.5
$ws.Range("K25").Value = 1651.5
$ws.Range("M25").Value = -1416.5

# Row 100: And My Axe / Doman Iron War Axe
$ws.Range("H100").Value = 33333
$ws.Range("J100").Value = 33333
$ws.Range("L100").Value = 33333
$ws.Range("N100").Value = -35497

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 1319.5264
$ws.Range("I107").Value = 1493.1538
$ws.Range("J107").Value = 943.3333
$ws.Range("K107").Value = 1493.1538
$ws.Range("L107").Value = 943.3333
$ws.Range("M107").Value = 426.8462
$ws.Range("N107").Value = -4783.3333

# Row 117: Idol Hands / Titanbronze Chakrams
$ws.Range("H117").Value = 58533.332
$ws.Range("J117").Value = 58533.332
$ws.Range("L117").Value = 58533.332
$ws.Range("N117").Value = -67711.33199999999

# Row 125: Archon of His Eye / High Durium Knives
$ws.Range("H125").Value = 79800
$ws.Range("J125").Value = 79800
$ws.Range("L125").Value = 79800
$ws.Range("N125").Value = -89640

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 17889.5
$ws.Range("I31").Value = 5555
$ws.Range("K31").Value = 5555
$ws.Range("M31").Value = -5260

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 17889.5
$ws.Range("I34").Value = 5555
$ws.Range("K34").Value = 5555
$ws.Range("M34").Value = -5353

# Row 94: Beech, Please / Beech Lumber
$ws.Range("H94").Value = 1398.9166
$ws.Range("I94").Value = 1624.6666
$ws.Range("J94").Value = 1323.6666
$ws.Range("K94").Value = 1624.6666
$ws.Range("L94").Value = 1323.6666
$ws.Range("M94").Value = -1173.6666
$ws.Range("N94").Value = -2225.6666

# Row 105: Zelkova, My Love / Zelkova Lumber
$ws.Range("H105").Value = 646
$ws.Range("I105").Value = 646
$ws.Range("K105").Value = 646
$ws.Range("M105").Value = 1101

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 2690.923
$ws.Range("I132").Value = 2503.7222
$ws.Range("J132").Value = 3112.125
$ws.Range("K132").Value = 7511.1666
$ws.Range("L132").Value = 9336.375
$ws.Range("M132").Value = -4981.1666
$ws.Range("N132").Value = -14396.375

$ws = $wb.Worksheets.Item("CUL")
# Row 119: Super Dark Times / Risotto al Nero
$ws.Range("H119").Value = 1437.8125
$ws.Range("I119").Value = 1336.8
$ws.Range("J119").Value = 1606.1666
$ws.Range("K119").Value = 4010.4
$ws.Range("L119").Value = 4818.4998
$ws.Range("M119").Value = 827.6000000000004
$ws.Range("N119").Value = -14494.4998

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 12965.468
$ws.Range("I131").Value = 479.54544
$ws.Range("J131").Value = 44180.273
$ws.Range("K131").Value = 1438.63632
$ws.Range("L131").Value = 132540.819
$ws.Range("M131").Value = 3601.36368
$ws.Range("N131").Value = -142620.819

# Row 132: More Mezcal / Cooking Mezcal
$ws.Range("H132").Value = 1613.9524
$ws.Range("I132").Value = 1457.7142
$ws.Range("J132").Value = 1926.4286
$ws.Range("K132").Value = 13119.4278
$ws.Range("L132").Value = 17337.8574
$ws.Range("M132").Value = -10589.4278
$ws.Range("N132").Value = -22397.8574

$ws = $wb.Worksheets.Item("LTW")
# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 4166.6665
$ws.Range("I93").Value = 3200
$ws.Range("K93").Value = 3200
$ws.Range("M93").Value = -1952

# Row 99: Shoe on the Other Foot / Tigerskin Boots of Crafting
$ws.Range("H99").Value = 51000
$ws.Range("J99").Value = 51000
$ws.Range("L99").Value = 51000
$ws.Range("N99").Value = -56990

$ws = $wb.Worksheets.Item("WVR")
# Row 81: Where the Dragonflies, the Net Catches / Crawler Silk
$ws.Range("H81").Value = 2905.4443
$ws.Range("I81").Value = 2562.375
$ws.Range("J81").Value = 5650
$ws.Range("K81").Value = 5124.75
$ws.Range("L81").Value = 11300
$ws.Range("M81").Value = -4063.75
$ws.Range("N81").Value = -13422

# Row 84: To Kill a Dragon on Nameday (L) / Crawler Silk
$ws.Range("H84").Value = 2905.4443
$ws.Range("I84").Value = 2562.375
$ws.Range("J84").Value = 5650
$ws.Range("K84").Value = 25623.75
$ws.Range("L84").Value = 56500
$ws.Range("M84").Value = -20319.75
$ws.Range("N84").Value = -67108

# Row 113: A Tender Table / Pixie Floss
$ws.Range("H113").Value = 6820
$ws.Range("I113").Value = 736.5
$ws.Range("K113").Value = 2209.5
$ws.Range("M113").Value = -39.5

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 2271
$ws.Range("I132").Value = 1409.3529
$ws.Range("K132").Value = 4228.0587
$ws.Range("M132").Value = -1698.0587

# Row 141: Silk for Sunperch / Thunderyards Silk Coat of Casting
$ws.Range("H141").Value = 47439.168
$ws.Range("J141").Value = 48965.453
$ws.Range("L141").Value = 48965.453
$ws.Range("N141").Value = -59325.453
